$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new daily log row (row 64): 2025/10/05, 日 (Sunday), hour 8, ranking 201.
# A leading apostrophe forces the date-like text to be stored as a literal string
# (matching the other "日付" column entries) instead of being auto-converted into a
# real date serial number. Resetting the style back to "Normal" afterwards keeps the
# new cell free of any explicit per-cell formatting, just like the surrounding rows.
$ws.Range("A64").Value = "'2025/10/05"
$ws.Range("A64").Style = "Normal"
$ws.Range("B64").Value = "日"
$ws.Range("C64").Value = 8
$ws.Range("D64").Value = 201
